# Re-style the three data tables (slides 14, 15, 16) from the deck's
# custom "Table_0" style to the built-in "No Style, No Grid" table
# style, matching the author's table-style-gallery selection.
#
# Table.ApplyStyle(StyleId) maps directly onto <a:tableStyleId> inside
# <a:tblPr> for the table's <a:tbl>.

$p = $ppt.ActivePresentation

$newStyleId = "{1B2168FB-3C16-4F38-B9B9-7B4829BD7B95}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
